$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "2025-03-15 13:19:21"
$ws.Range("E2").Value = "POST"
$ws.Range("F2").Value = "http://49.234.6.241:5230/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("G2").Value = "/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("M2").Value = 0.004
$ws.Range("N2").Value = 0
$ws.Range("Q2").Value = $true

# Row 3
$ws.Range("D3").Value = "2025-03-15 13:19:21"
$ws.Range("E3").Value = "POST"
$ws.Range("F3").Value = "http://49.234.6.241:5230/memos.api.v1.MemoService/CreateMemo"
$ws.Range("G3").Value = "/memos.api.v1.MemoService/CreateMemo"
$ws.Range("M3").Value = 0.003
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $true
$ws.Range("Q3").Value = $true

# Row 4
$ws.Range("D4").Value = "2025-03-15 13:19:21"
$ws.Range("E4").Value = "POST"
$ws.Range("F4").Value = "http://49.234.6.241:5230/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("G4").Value = "/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("M4").Value = 0.002

# Row 5
$ws.Range("D5").Value = "2025-03-15 13:19:21"
$ws.Range("E5").Value = "POST"
$ws.Range("F5").Value = "http://49.234.6.241:5230/memos.api.v1.MemoService/CreateMemo"
$ws.Range("G5").Value = "/memos.api.v1.MemoService/CreateMemo"
$ws.Range("M5").Value = 0.002

# Row 6
$ws.Range("D6").Value = "2025-03-15 13:19:21"
$ws.Range("E6").Value = "POST"
$ws.Range("F6").Value = "http://49.234.6.241:5230/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("G6").Value = "/memos.api.v1.UserService/CreateUserAccessToken"

# Row 7
$ws.Range("D7").Value = "2025-03-15 13:19:21"
$ws.Range("E7").Value = "POST"
$ws.Range("F7").Value = "http://49.234.6.241:5230/memos.api.v1.MemoService/CreateMemo"
$ws.Range("G7").Value = "/memos.api.v1.MemoService/CreateMemo"
$ws.Range("M7").Value = 0.002

# Row 8
$ws.Range("D8").Value = "2025-03-15 13:19:21"
$ws.Range("E8").Value = "POST"
$ws.Range("F8").Value = "http://49.234.6.241:5230/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("G8").Value = "/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("M8").Value = 0.002
$ws.Range("N8").Value = 1
$ws.Range("Q8").Value = $false

# Row 9
$ws.Range("D9").Value = "2025-03-15 13:19:21"
$ws.Range("E9").Value = "POST"
$ws.Range("F9").Value = "http://49.234.6.241:5230/memos.api.v1.MemoService/CreateMemo"
$ws.Range("G9").Value = "/memos.api.v1.MemoService/CreateMemo"
$ws.Range("M9").Value = 0.003
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = $false
